$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 876, shifting rows 876:917 down to 877:918
$ws.Rows.Item(876).Insert()

# Populate the newly inserted row 876 with the new data point
# (date written with a leading apostrophe so Excel stores it as literal
# text like the rest of the column instead of auto-converting it to a
# date serial; ClearFormats keeps the cell formatted as "General" instead
# of leaving the "text" quote-prefix formatting behind)
$ws.Range("A876").Value = "'2026/02/24"
$ws.Range("A876").ClearFormats()
$ws.Range("B876").Value = "火"
$ws.Range("C876").Value = 20
$ws.Range("D876").Value = 201
